$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rotate the pre-existing per-row cell styles in C2:D10 by one slot ---
# (Mirrors how this log's formatting drifted historically: inserting a row
#  pushes the existing style objects for rows 2-10 down into rows 3-11, and
#  removing the now-redundant trailing row pulls row 11 back out, leaving the
#  same 9 pre-existing styles occupying rows 2-10 but rotated by one.)
$ws.Rows("2:2").Insert()
$ws.Rows("11:11").Delete()

# --- Step 2: the row insert/delete also shifted the cell VALUES down by one
#     row; restore each existing entry's data to its original row (Value
#     assignment alone does not touch the style that now sits in that row). ---
$ws.Range("A2").Value = (Get-Date -Year 2017 -Month 4 -Day 17 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B2").Value = 1.5
$ws.Range("C2").Value = "UI"
$ws.Range("D2").Value = "Fixed some of the stuff per Prof Sedlemeyer's comments on UI"

$ws.Range("A3").Value = (Get-Date -Year 2017 -Month 4 -Day 18 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = "TopTen.txt"
$ws.Range("D3").Value = "Changed the way we save top scores to a text file"

$ws.Range("A4").Value = (Get-Date -Year 2017 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "SF-17"
$ws.Range("D4").Value = "Almost finished everything for this user story. All that is left is a sort method for the model"

$ws.Range("A5").Value = (Get-Date -Year 2017 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "SF-17"
$ws.Range("D5").Value = "Created the sort method and properly implemented it"

$ws.Range("A6").Value = (Get-Date -Year 2017 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B6").Value = 1.5
$ws.Range("C6").Value = "SF-17"
$ws.Range("D6").Value = "Added the date to both top ten lists"

$ws.Range("A7").Value = (Get-Date -Year 2017 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "SF-13"
$ws.Range("D7").Value = "Added a button that will remove all tiles of a number based on user input"

# --- Step 3: fill in the new Sprint entry for row 8 (previously blank) ---
$ws.Range("A8").Value = (Get-Date -Year 2017 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "SF-13 & SF-14"
$ws.Range("D8").Value = "Made the changes for the limit on how often you can use these features"

# --- Step 4: rows 9 and 10 go back to being blank (their old data, if any,
#     was only the leftover shifted duplicate from row 8/9) ---
$ws.Range("A9:D9").ClearContents()
$ws.Range("A10:D10").ClearContents()

# Update the active selection to reflect where the author finished editing
$ws.Range("C9").Select()
